$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").Value = 20220401
$ws.Range("B76").Value = 2225.92
$ws.Range("C76").Value = 2224.4699999999998
$ws.Range("E76").Value = 180
$ws.Range("F76").Value = "CRM OPENED 20220401 MG"

$ws.Range("A77").Value = 20220404
$ws.Range("B77").Value = 2224.5797358742998
$ws.Range("C77").Value = 2224.4699999999998
$ws.Range("E77").Value = 180
$ws.Range("F77").Value = "CRM OPENED 20220401 MG"

$ws.Range("D70:D77").Formula = "=100*(B70-C70)/C70"

$ws.Range("D77").Select()
